$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null

$ws.Range("A10").Value = "Sphagetti"
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = "NTU"
$ws.Range("D10").Value = "pasta"

$ws.Range("A11").Value = "cheeseburger"
$ws.Range("B11").Value = 1.5
$ws.Range("C11").Value = "NTU"
$ws.Range("D11").Value = "burger"

$ws.Range("A10:D11").Select()
